# Bug 487 / #534: the HTML-to-Word generator was emitting a spurious
# trailing manual line break (<w:br/>) at the end of the paragraph that
# introduces each bulleted list ("See the following list:" and
# "A new list of:"). The fix ignores/removes that extra <br> so the
# paragraph ends right after the introductory text, instead of carrying
# an empty trailing run with a line break into the following list.
#
# ^l in a Word Find pattern matches a manual line break character
# (the <w:br/> run), so searching for "<text>^l" and replacing with
# "<text>" deletes that trailing run while leaving everything else
# (including the following bulleted-list paragraphs) untouched.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "See the following list:^l", $true, $false, $false, $false, $false,
    $true, 1, $false, "See the following list:", 2)

$d.Content.Find.Execute(
    "A new list of:^l", $true, $false, $false, $false, $false,
    $true, 1, $false, "A new list of:", 2)
